# Update the cryptos list (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (matches the workbook's inlineStr cells),
# preventing Excel from auto-coercing numeric-looking strings (e.g. "0.497")
# into actual numbers, while leaving the cell's style untouched afterwards.
function Set-Text($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

function Set-Pct($row, $text) {
    Set-Text $row 5 "  $text  "
}

# Row 2 - Bitcoin
Set-Text 2 4 "26.196.29"
Set-Pct 2 "-1.97%"

# Row 3 - Ethereum
Set-Text 3 4 "1.582.45"
Set-Pct 3 "-1.17%"

# Row 4 - TetherUSD
Set-Pct 4 "-0.37%"

# Row 5 - BNB
Set-Text 5 4 "209.65"
Set-Pct 5 "-1.03%"

# Row 6 - XRP
Set-Text 6 4 "0.497"
Set-Pct 6 "-2.91%"

# Row 7 - USDC
Set-Pct 7 "-0.34%"

# Row 8 - Dogecoin
Set-Pct 8 "-1.39%"

# Row 9 - Cardano
Set-Text 9 4 "0.246"
Set-Pct 9 "-0.59%"

# Row 10 - Solana
Set-Text 10 4 "19.51"
Set-Pct 10 "-1.12%"

# Row 11 - TRON
Set-Text 11 4 "0.0846"
Set-Pct 11 "+0.14%"

# Row 12 - WrappedliquidstakedEther2.0
Set-Text 12 4 "1.804.56"
Set-Pct 12 "-1.21%"

# Row 13 - WrappedEther
Set-Text 13 4 "1.577.27"
Set-Pct 13 "-1.25%"

# Row 14 - Polkadot
Set-Text 14 4 "4.04"
Set-Pct 14 "+0.07%"

# Row 15 - Polygon
Set-Pct 15 "-1.39%"

# Row 16 - Litecoin
Set-Text 16 4 "64.49"
Set-Pct 16 "-0.67%"

# Row 17 - WrappedBTC
Set-Text 17 4 "26.197.43"
Set-Pct 17 "-1.85%"

# Row 18 - ShibaInu
Set-Text 18 4 "0.0₃0734"
Set-Pct 18 "-0.89%"

# Row 19 - Chainlink
Set-Text 19 4 "7.26"
Set-Pct 19 "+1.10%"

# Row 20 - Dai
Set-Pct 20 "-0.34%"

# Row 21 - BitcoinCash
Set-Text 21 4 "206.82"
Set-Pct 21 "-1.73%"

# Row 22 - Uniswap
Set-Pct 22 "-0.57%"

# Row 23 - Toncoin
Set-Pct 23 "-3.17%"

# Row 24 - Avalanche
Set-Text 24 4 "8.87"
Set-Pct 24 "-0.96%"

# Row 25 - Monero
Set-Text 25 4 "144.90"
Set-Pct 25 "+0.60%"

# Row 26 - BinanceUSD
Set-Pct 26 "-0.24%"

# Row 27 - Cosmos
Set-Text 27 4 "7.03"
Set-Pct 27 "-0.62%"

# Row 28 - Stellar
Set-Pct 28 "-1.06%"

# Row 29 - EthereumClassic
Set-Text 29 4 "15.20"

# Row 30 - Hedera
Set-Pct 30 "-1.24%"

# Row 31 - PancakeSwap
Set-Text 31 4 "1.14"
Set-Pct 31 "-1.00%"

# Row 32 - Filecoin
Set-Pct 32 "-1.32%"

# Row 33 - InternetComputer(DFINITY)
Set-Text 33 4 "2.95"
Set-Pct 33 "-0.88%"

# Row 34 - Maker
Set-Text 34 4 "1.282.14"
Set-Pct 34 "-0.75%"

# Row 35 - now HuobiToken (was WEMIXToken)
Set-Text 35 2 "HuobiToken"
Set-Text 35 3 "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-Text 35 4 "2.46"
Set-Pct 35 "-0.42%"

# Row 36 - now WEMIXToken (was HuobiToken)
Set-Text 36 2 "WEMIXToken"
Set-Text 36 3 "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-Text 36 4 "1.22"
Set-Pct 36 "+6.89%"

# Row 37 - ImmutableX
Set-Text 37 4 "0.608"
Set-Pct 37 "+1.58%"

# Row 38 - LidoDAOToken
Set-Pct 38 "-0.87%"

# Row 39 - VeChain
Set-Pct 39 "-1.39%"

# Row 40 - ARBITRUM
Set-Text 40 4 "0.815"
Set-Pct 40 "-1.72%"

# Row 41 - FraxShare
Set-Text 41 4 "5.58"
Set-Pct 41 "+3.36%"

# Row 42 - TrustWalletToken
Set-Text 42 4 "0.767"
Set-Pct 42 "-1.88%"

# Row 43 - MXToken
Set-Pct 43 "-3.02%"

# Row 44 - Aave
Set-Text 44 4 "62.27"
Set-Pct 44 "-1.20%"

# Row 45 - RocketPoolETH
Set-Text 45 4 "1.718.35"
Set-Pct 45 "-1.23%"

# Row 46 - Quant
Set-Text 46 4 "88.84"
Set-Pct 46 "-1.95%"

# Row 47 - RenderToken
Set-Pct 47 "-0.17%"

# Row 48 - Algorand
Set-Pct 48 "-0.89%"

# Row 49 - Cronos
Set-Text 49 4 "0.0507"
Set-Pct 49 "-1.70%"

# Row 50 - BabyDogeCoin
Set-Text 50 4 "0.0₇0957"
Set-Pct 50 "-10.10%"

# Row 51 - USDD
Set-Pct 51 "-0.22%"
